$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on C5 so "67%" is stored as the literal string
# "67%" (matching the source data) instead of being auto-converted to a
# numeric percentage value by Excel's type inference.
$ws.Range("C5").NumberFormat = "@"

$ws.Range("A5").Value = "2025-07-01 12:20:32"
$ws.Range("B5").Value = "16°C"
$ws.Range("C5").Value = "67%"
